$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 02:44"

# --- Update country case numbers (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7991999
$ws.Range("C4").Value = 41936
$ws.Range("D4").Value = 5128162
$ws.Range("E4").Value = 2644142
$ws.Range("G4").Value = 325
$ws.Range("H4").Value = 219695

# Row 6: Brasil
$ws.Range("B6").Value = 5094979
$ws.Range("C6").Value = 3139
$ws.Range("D6").Value = 4470165
$ws.Range("E6").Value = 474308
$ws.Range("G6").Value = 270
$ws.Range("H6").Value = 150506

# Row 11: Peru
$ws.Range("B11").Value = 849371
$ws.Range("C11").Value = 3283
$ws.Range("D11").Value = 743969
$ws.Range("E11").Value = 72097
$ws.Range("G11").Value = 82
$ws.Range("H11").Value = 33305

# Row 25: Alemania
$ws.Range("D25").Value = 274700
$ws.Range("E25").Value = 41889

# Row 117: Mauritania
$ws.Range("D117").Value = 7274
$ws.Range("E117").Value = 113

# Rows 129-132: Bahamas, Guinea Ecuatorial, Surinam, Trinidad yTobago were
# reordered (Bahamas now before Guinea Ecuatorial, Surinam before Trinidad
# yTobago) and their data updated.
$ws.Range("A129").Value = "Bahamas"
$ws.Range("B129").Value = 5078
$ws.Range("C129").Value = 123
$ws.Range("D129").Value = 2900
$ws.Range("E129").Value = 2071
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 107

$ws.Range("A130").Value = "Guinea Ecuatorial"
$ws.Range("B130").Value = 5063
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 4894
$ws.Range("E130").Value = 86
$ws.Range("H130").Value = 83

$ws.Range("A131").Value = "Surinam"
$ws.Range("B131").Value = 5051
$ws.Range("C131").Value = 16
$ws.Range("D131").Value = 4845
$ws.Range("E131").Value = 99

$ws.Range("A132").Value = "Trinidad yTobago"
$ws.Range("B132").Value = 5043
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 3221
$ws.Range("E132").Value = 1732
$ws.Range("H132").Value = 90

# Row 141: Somalia
$ws.Range("B141").Value = 3864
$ws.Range("C141").Value = 17
$ws.Range("D141").Value = 3089
$ws.Range("E141").Value = 676

# Row 169: Santo Tome y Principe
$ws.Range("B169").Value = 929
$ws.Range("C169").Value = 7
$ws.Range("D169").Value = 892
$ws.Range("E169").Value = 22

# Row 185: Isla de Man
$ws.Range("B185").Value = 346
$ws.Range("C185").Value = 1
$ws.Range("E185").Value = 6

# Row 191: Barbados
$ws.Range("B191").Value = 208
$ws.Range("C191").Value = 2
$ws.Range("D191").Value = 186
$ws.Range("E191").Value = 15
